$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 10842.134
$ws.Range("J32").Value = 14773
$ws.Range("L32").Value = 14773
$ws.Range("N32").Value = -15425
$ws.Range("H40").Value = 2549.5
$ws.Range("J40").Value = 2924.25
$ws.Range("L40").Value = 2924.25
$ws.Range("N40").Value = -3274.25
$ws.Range("H41").Value = 686.25714
$ws.Range("I41").Value = 236.66667
$ws.Range("K41").Value = 236.66667
$ws.Range("M41").Value = 203.33333
$ws.Range("H76").Value = 66673840
$ws.Range("I76").Value = 142864820
$ws.Range("J76").Value = 6740.5
$ws.Range("K76").Value = 142864820
$ws.Range("L76").Value = 6740.5
$ws.Range("M76").Value = -142864505
$ws.Range("N76").Value = -7370.5
$ws.Range("H79").Value = 66673840
$ws.Range("I79").Value = 142864820
$ws.Range("J79").Value = 6740.5
$ws.Range("K79").Value = 142864820
$ws.Range("L79").Value = 6740.5
$ws.Range("M79").Value = -142863728
$ws.Range("N79").Value = -8924.5
$ws.Range("H86").Value = 125004824
$ws.Range("I86").Value = 250001780
$ws.Range("J86").Value = 7874.25
$ws.Range("K86").Value = 250001780
$ws.Range("L86").Value = 7874.25
$ws.Range("M86").Value = -250000657
$ws.Range("N86").Value = -10120.25
$ws.Range("H88").Value = 2633.4
$ws.Range("I88").Value = 2466.8333
$ws.Range("J88").Value = 2704.7856
$ws.Range("K88").Value = 2466.8333
$ws.Range("L88").Value = 2704.7856
$ws.Range("M88").Value = -2060.8333
$ws.Range("N88").Value = -3516.7856
$ws.Range("H89").Value = 125004824
$ws.Range("I89").Value = 250001780
$ws.Range("J89").Value = 7874.25
$ws.Range("K89").Value = 1250008900
$ws.Range("L89").Value = 39371.25
$ws.Range("M89").Value = -1250003284
$ws.Range("N89").Value = -50603.25
$ws.Range("H91").Value = 2633.4
$ws.Range("I91").Value = 2466.8333
$ws.Range("J91").Value = 2704.7856
$ws.Range("K91").Value = 2466.8333
$ws.Range("L91").Value = 2704.7856
$ws.Range("M91").Value = -1062.8333
$ws.Range("N91").Value = -5512.7856
$ws.Range("H112").Value = 2402.3044
$ws.Range("J112").Value = 2555.5264
$ws.Range("L112").Value = 7666.5792
$ws.Range("N112").Value = -9882.5792
$ws.Range("H116").Value = 6177387.5
$ws.Range("I116").Value = 9261082
$ws.Range("K116").Value = 9261082
$ws.Range("M116").Value = -9257640
$ws.Range("H138").Value = 6822.5947
$ws.Range("I138").Value = 4383.8887
$ws.Range("J138").Value = 7606.4644
$ws.Range("K138").Value = 13151.6661
$ws.Range("L138").Value = 22819.3932
$ws.Range("M138").Value = -8011.666100000002
$ws.Range("N138").Value = -33099.3932

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18560180
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 18560180
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 18560180
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -18560754
$ws.Range("H60").Value = 41995
$ws.Range("J60").Value = 41995
$ws.Range("L60").Value = 41995
$ws.Range("N60").Value = -43461
$ws.Range("H63").Value = 4741.25
$ws.Range("I63").Value = 3275.7144
$ws.Range("J63").Value = 15000
$ws.Range("K63").Value = 3275.7144
$ws.Range("L63").Value = 15000
$ws.Range("M63").Value = -2589.7144
$ws.Range("N63").Value = -16372
$ws.Range("H66").Value = 4741.25
$ws.Range("I66").Value = 3275.7144
$ws.Range("J66").Value = 15000
$ws.Range("K66").Value = 16378.572
$ws.Range("L66").Value = 75000
$ws.Range("M66").Value = -12946.572
$ws.Range("N66").Value = -81864
$ws.Range("H88").Value = 6312.1665
$ws.Range("I88").Value = 4719.2
$ws.Range("J88").Value = 7450
$ws.Range("K88").Value = 4719.2
$ws.Range("L88").Value = 7450
$ws.Range("M88").Value = -4313.2
$ws.Range("N88").Value = -8262
$ws.Range("H91").Value = 6312.1665
$ws.Range("I91").Value = 4719.2
$ws.Range("J91").Value = 7450
$ws.Range("K91").Value = 4719.2
$ws.Range("L91").Value = 7450
$ws.Range("M91").Value = -3315.2
$ws.Range("N91").Value = -10258
$ws.Range("H102").Value = 3472.8572
$ws.Range("I102").Value = 3571.5
$ws.Range("J102").Value = 1500
$ws.Range("K102").Value = 3571.5
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = -1949.5
$ws.Range("N102").Value = -4744
$ws.Range("H132").Value = 707434.7
$ws.Range("I132").Value = 861543.6
$ws.Range("J132").Value = 79144.38
$ws.Range("K132").Value = 2584630.8
$ws.Range("L132").Value = 237433.14
$ws.Range("M132").Value = -2582100.8
$ws.Range("N132").Value = -242493.14

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 11729.857
$ws.Range("I86").Value = 1485
$ws.Range("K86").Value = 1485
$ws.Range("M86").Value = -362
$ws.Range("H89").Value = 11729.857
$ws.Range("I89").Value = 1485
$ws.Range("K89").Value = 7425
$ws.Range("M89").Value = -1809
$ws.Range("H134").Value = 1037921.8
$ws.Range("I134").Value = 1136411.5
$ws.Range("J134").Value = 13629.6
$ws.Range("K134").Value = 3409234.5
$ws.Range("L134").Value = 40888.8
$ws.Range("M134").Value = -3406699.5
$ws.Range("N134").Value = -45958.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 25540.834
$ws.Range("I31").Value = 28712.857
$ws.Range("J31").Value = 21100
$ws.Range("K31").Value = 28712.857
$ws.Range("L31").Value = 21100
$ws.Range("M31").Value = -28417.857
$ws.Range("N31").Value = -21690
$ws.Range("H34").Value = 25540.834
$ws.Range("I34").Value = 28712.857
$ws.Range("J34").Value = 21100
$ws.Range("K34").Value = 28712.857
$ws.Range("L34").Value = 21100
$ws.Range("M34").Value = -28510.857
$ws.Range("N34").Value = -21504
$ws.Range("H62").Value = 11597.833
$ws.Range("I62").Value = 13529.333
$ws.Range("J62").Value = 9666.333000000001
$ws.Range("K62").Value = 13529.333
$ws.Range("L62").Value = 9666.333000000001
$ws.Range("M62").Value = -12905.333
$ws.Range("N62").Value = -10914.333
$ws.Range("H65").Value = 11597.833
$ws.Range("I65").Value = 13529.333
$ws.Range("J65").Value = 9666.333000000001
$ws.Range("K65").Value = 67646.66500000001
$ws.Range("L65").Value = 48331.665
$ws.Range("M65").Value = -64526.66500000001
$ws.Range("N65").Value = -54571.665
$ws.Range("H132").Value = 4061.457
$ws.Range("I132").Value = 2986.3872
$ws.Range("K132").Value = 8959.161599999999
$ws.Range("M132").Value = -6429.161599999999
$ws.Range("H141").Value = 405606.34
$ws.Range("J141").Value = 405606.34
$ws.Range("L141").Value = 405606.34
$ws.Range("N141").Value = -415966.34

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1583.72
$ws.Range("I5").Value = 780.0769
$ws.Range("J5").Value = 2454.3333
$ws.Range("K5").Value = 2340.2307
$ws.Range("L5").Value = 7362.999899999999
$ws.Range("M5").Value = -2228.2307
$ws.Range("N5").Value = -7586.999899999999
$ws.Range("H94").Value = 11754.272
$ws.Range("I94").Value = 1999.5
$ws.Range("J94").Value = 13922
$ws.Range("K94").Value = 5998.5
$ws.Range("L94").Value = 41766
$ws.Range("M94").Value = -5322.5
$ws.Range("N94").Value = -43118
$ws.Range("H110").Value = 25949.75
$ws.Range("I110").Value = 13799
$ws.Range("K110").Value = 41397
$ws.Range("M110").Value = -37307
$ws.Range("H111").Value = 5934.25
$ws.Range("I111").Value = 3210.5715
$ws.Range("K111").Value = 9631.7145
$ws.Range("M111").Value = -6564.7145
$ws.Range("H112").Value = 12540.75
$ws.Range("I112").Value = 630.6667
$ws.Range("J112").Value = 16510.777
$ws.Range("K112").Value = 1892.0001
$ws.Range("L112").Value = 49532.33099999999
$ws.Range("M112").Value = -784.0001
$ws.Range("N112").Value = -51748.33099999999
$ws.Range("H135").Value = 1583.72
$ws.Range("I135").Value = 780.0769
$ws.Range("J135").Value = 2454.3333
$ws.Range("K135").Value = 7020.6921
$ws.Range("L135").Value = 22088.9997
$ws.Range("M135").Value = -4485.6921
$ws.Range("N135").Value = -27158.9997
$ws.Range("H136").Value = 18521906
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 71437416
$ws.Range("I80").Value = 200004590
$ws.Range("J80").Value = 11207.667
$ws.Range("K80").Value = 200004590
$ws.Range("L80").Value = 11207.667
$ws.Range("M80").Value = -200003592
$ws.Range("N80").Value = -13203.667
$ws.Range("H83").Value = 71437416
$ws.Range("I83").Value = 200004590
$ws.Range("J83").Value = 11207.667
$ws.Range("K83").Value = 1000022950
$ws.Range("L83").Value = 56038.335
$ws.Range("M83").Value = -1000017958
$ws.Range("N83").Value = -66022.33499999999
$ws.Range("H132").Value = 6336.098
$ws.Range("I132").Value = 5752.386
$ws.Range("K132").Value = 17257.158
$ws.Range("M132").Value = -14727.158

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6095.154
$ws.Range("I40").Value = 5862.9546
$ws.Range("J40").Value = 7372.25
$ws.Range("K40").Value = 5862.9546
$ws.Range("L40").Value = 7372.25
$ws.Range("M40").Value = -5726.9546
$ws.Range("N40").Value = -7644.25
$ws.Range("H46").Value = 45455908
$ws.Range("I46").Value = 999
$ws.Range("J46").Value = 83335000
$ws.Range("K46").Value = 999
$ws.Range("L46").Value = 83335000
$ws.Range("M46").Value = -811
$ws.Range("N46").Value = -83335376
$ws.Range("H136").Value = 8795.1
$ws.Range("I136").Value = 9440.817999999999
$ws.Range("J136").Value = 8421.263000000001
$ws.Range("K136").Value = 28322.454
$ws.Range("L136").Value = 25263.789
$ws.Range("M136").Value = -25772.454
$ws.Range("N136").Value = -30363.789

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 9202.983
$ws.Range("I122").Value = 5229.878
$ws.Range("K122").Value = 15689.634
$ws.Range("M122").Value = -13239.634
$ws.Range("H132").Value = 8070.1577
$ws.Range("I132").Value = 6820.722
$ws.Range("K132").Value = 20462.166
$ws.Range("M132").Value = -17932.166
$ws.Range("H136").Value = 8781952
$ws.Range("I136").Value = 11372128
$ws.Range("J136").Value = 15202.538
$ws.Range("K136").Value = 34116384
$ws.Range("L136").Value = 45607.614
$ws.Range("M136").Value = -34113834
$ws.Range("N136").Value = -50707.614

Write-Output "done"